$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C10 (row 10, "R30" rule's "From" value) from 18 to 1
$ws.Range("C10").Value = 1
